# Update the "Init" sheet's lower-right cell references (column D, rows 5-11)
# so that they point at row 36 + 3 = row 39 instead of row 36.
# This reflects additional scenario rows being added to the underlying
# listing tables referenced by this init file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value  = "A39"
$ws.Range("D6").Value  = "B39"
$ws.Range("D7").Value  = "C39"
$ws.Range("D8").Value  = "G39"
$ws.Range("D9").Value  = "H39"
$ws.Range("D10").Value = "I39"
$ws.Range("D11").Value = "J39"
